# Hindalco prices worksheet update (2025-08-26 06:50:27 UTC)
# Inserts a new latest-price row (Sl.no. 27, 26.08.2025 circular) above the
# existing data, shifting every other row down by one and re-pointing the
# hyperlinks in column F so that row 2..row 8 each keep their circular link.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel constants
$xlPasteFormats = -4122

# All existing hyperlinks are anchored to fixed cell refs (F2..F7) and won't
# move on their own when we insert a row, so drop them now - we'll recreate
# the full, correctly-shifted set at the end.
$ws.Hyperlinks.Delete()

# Push every data row (old rows 2-27) down by one to make room for the new
# latest entry, then restore the data-row look (non-bold, centered) and
# number format on the newly inserted row 2 by copying formatting from the
# row right below it (old row 2, now row 3).
$ws.Rows.Item(2).Insert()
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# Fill in the brand-new row 2 with the latest circular figures.
$ws.Range("A2").Value2 = 27
$ws.Range("B2").Value2 = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Range("C2").Value2 = "P1020"
$ws.Range("D2").Value2 = 258.25
$ws.Range("E2").Value2 = "26.08.2025"
$ws.Range("F2").Value2 = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-26-august-2025.pdf"

# Dimension should now cover the extra row.
Write-Output ("UsedRange=" + $ws.UsedRange.Address())

# Re-create the circular-link hyperlinks for rows 2-8 (the new row plus the
# six rows that used to be 2-7), in order, so link ids line up with the rows.
$links = @(
    @{ Row = 2; Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-26-august-2025.pdf" },
    @{ Row = 3; Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-23-august-2025.pdf" },
    @{ Row = 4; Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-20-august-2025.pdf" },
    @{ Row = 5; Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-19-august-2025.pdf" },
    @{ Row = 6; Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-14-august-2025.pdf" },
    @{ Row = 7; Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-13-august-2025.pdf" },
    @{ Row = 8; Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-august-2025.pdf" }
)

# A reference cell that already carries the plain (non-hyperlink) look used
# throughout column F, so we can restore it after Hyperlinks.Add recolors
# the cell with the built-in "Hyperlink" style.
$plainFormatCell = $ws.Range("E2")

foreach ($link in $links) {
    $cell = $ws.Cells.Item($link.Row, 6)
    $ws.Hyperlinks.Add($cell, $link.Url)
    $plainFormatCell.Copy()
    $cell.PasteSpecial($xlPasteFormats)
}
$excel.CutCopyMode = 0

Write-Output "Hindalco sheet updated with 26.08.2025 circular row."
